$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: move the "polity2" row (row 5, under Statehood) to the bottom of the table
# and relabel it with its own Concept Label "Polity"
$ws.Range("A19:E19").Insert(-4121)
$ws.Range("A5:E5").Cut($ws.Range("A19:E19"))
$ws.Range("A5:E5").Delete(-4162)
$ws.Cells.Item(18,1).Value = "Polity"

# Step 2: insert 2 new rows between "FD / Flood Dummy" (row 6) and "DD / Drought Dummy" (row 7)
# for the new "SD / Storm Dummy" and "ED / Earthquake Dummy" concepts
$ws.Rows.Item(7).Resize(2).Insert()
$ws.Cells.Item(7,2).Value = "SD"
$ws.Cells.Item(7,3).Value = "Storm Dummy"
$ws.Cells.Item(7,5).Value = "http://www.emdat.be/database"
$ws.Cells.Item(8,2).Value = "ED"
$ws.Cells.Item(8,3).Value = "Earthquake Dummy"
$ws.Cells.Item(8,5).Value = "http://www.emdat.be/database"

# Step 3: update the selected cell to match the final edit position
$ws.Range("E8").Select()

Write-Output "done"
